# Dev Diary: add a new "13/06/18" entry after the existing last entry,
# and move the "_GoBack" bookmark onto the new bullet's text.

$d = $word.ActiveDocument

# 1. Drop the existing (hidden) _GoBack bookmark from the end of the
#    "... int array" paragraph - it will be re-created around the new
#    text that becomes the new end of the document.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. Append the two new paragraphs (a Heading2 date heading followed by
#    a bulleted ListParagraph entry) at the very end of the document, so
#    the existing content is left untouched.
$endOfDoc = $d.Content.End
$insertionPoint = $d.Range($endOfDoc, $endOfDoc)

$newEntryXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading2"/>
            </w:pPr>
            <w:r>
              <w:t>13/06/18</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Added visualisation of the map that is stored</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($newEntryXml)

# 3. Re-create the _GoBack bookmark around the newly-inserted bullet text,
#    mirroring its original placement at the end of the document content.
$newLastParagraph = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $newLastParagraph.Range)
